$d = $word.ActiveDocument

$ids = @("p165r_1", "p165r_2", "p165r_3", "p165r_4", "p165r_5", "p165r_6")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $old, 2)
}
